$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "A twin jet narrowbody aircraft manufactured by Airbus"
$ws.Range("E11").Select()
